$d = $word.ActiveDocument

# Update the date heading (wdReplaceOne = 1 so we touch only this one run)
$d.Content.Find.Execute("2023-10-21 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-22 Sunday", 1) | Out-Null

# Update each arithmetic answer cell in the 20x5 table. The Find is
# scoped to each cell's own Range, and Replace=1 (wdReplaceOne) is used
# so that cells which happen to share identical old text (e.g. the two
# "48-19=29" cells) are each only replaced with their own new value.
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1)
$c.Range.Find.Execute("90-65=25", $true, $false, $false, $false, $false, $true, 1, $false, "27+44=71", 1) | Out-Null
$c = $t.Cell(1, 2)
$c.Range.Find.Execute("31+6=37", $true, $false, $false, $false, $false, $true, 1, $false, "37+2=39", 1) | Out-Null
$c = $t.Cell(1, 3)
$c.Range.Find.Execute("17+44=61", $true, $false, $false, $false, $false, $true, 1, $false, "64-3=61", 1) | Out-Null
$c = $t.Cell(1, 4)
$c.Range.Find.Execute("61+33=94", $true, $false, $false, $false, $false, $true, 1, $false, "5+11=16", 1) | Out-Null
$c = $t.Cell(1, 5)
$c.Range.Find.Execute("41-4=37", $true, $false, $false, $false, $false, $true, 1, $false, "45-35=10", 1) | Out-Null
$c = $t.Cell(2, 1)
$c.Range.Find.Execute("83-71=12", $true, $false, $false, $false, $false, $true, 1, $false, "46+51=97", 1) | Out-Null
$c = $t.Cell(2, 2)
$c.Range.Find.Execute("82+0=82", $true, $false, $false, $false, $false, $true, 1, $false, "43+20=63", 1) | Out-Null
$c = $t.Cell(2, 3)
$c.Range.Find.Execute("27+46=73", $true, $false, $false, $false, $false, $true, 1, $false, "91-42=49", 1) | Out-Null
$c = $t.Cell(2, 4)
$c.Range.Find.Execute("85-36=49", $true, $false, $false, $false, $false, $true, 1, $false, "52+33=85", 1) | Out-Null
$c = $t.Cell(2, 5)
$c.Range.Find.Execute("68-44=24", $true, $false, $false, $false, $false, $true, 1, $false, "59+35=94", 1) | Out-Null
$c = $t.Cell(3, 1)
$c.Range.Find.Execute("27-7=20", $true, $false, $false, $false, $false, $true, 1, $false, "25+13=38", 1) | Out-Null
$c = $t.Cell(3, 2)
$c.Range.Find.Execute("35+59=94", $true, $false, $false, $false, $false, $true, 1, $false, "69+20=89", 1) | Out-Null
$c = $t.Cell(3, 3)
$c.Range.Find.Execute("43-12=31", $true, $false, $false, $false, $false, $true, 1, $false, "18-15=3", 1) | Out-Null
$c = $t.Cell(3, 4)
$c.Range.Find.Execute("25+51=76", $true, $false, $false, $false, $false, $true, 1, $false, "21-19=2", 1) | Out-Null
$c = $t.Cell(3, 5)
$c.Range.Find.Execute("88-74=14", $true, $false, $false, $false, $false, $true, 1, $false, "74-12=62", 1) | Out-Null
$c = $t.Cell(4, 1)
$c.Range.Find.Execute("91+0=91", $true, $false, $false, $false, $false, $true, 1, $false, "76-31=45", 1) | Out-Null
$c = $t.Cell(4, 2)
$c.Range.Find.Execute("41+13=54", $true, $false, $false, $false, $false, $true, 1, $false, "13+85=98", 1) | Out-Null
$c = $t.Cell(4, 3)
$c.Range.Find.Execute("82-59=23", $true, $false, $false, $false, $false, $true, 1, $false, "27+51=78", 1) | Out-Null
$c = $t.Cell(4, 4)
$c.Range.Find.Execute("92-87=5", $true, $false, $false, $false, $false, $true, 1, $false, "5+24=29", 1) | Out-Null
$c = $t.Cell(4, 5)
$c.Range.Find.Execute("89-72=17", $true, $false, $false, $false, $false, $true, 1, $false, "1+45=46", 1) | Out-Null
$c = $t.Cell(5, 1)
$c.Range.Find.Execute("75+12=87", $true, $false, $false, $false, $false, $true, 1, $false, "24+1=25", 1) | Out-Null
$c = $t.Cell(5, 2)
$c.Range.Find.Execute("39-4=35", $true, $false, $false, $false, $false, $true, 1, $false, "16+82=98", 1) | Out-Null
$c = $t.Cell(5, 3)
$c.Range.Find.Execute("84-49=35", $true, $false, $false, $false, $false, $true, 1, $false, "6+42=48", 1) | Out-Null
$c = $t.Cell(5, 4)
$c.Range.Find.Execute("98-38=60", $true, $false, $false, $false, $false, $true, 1, $false, "83-38=45", 1) | Out-Null
$c = $t.Cell(5, 5)
$c.Range.Find.Execute("32-13=19", $true, $false, $false, $false, $false, $true, 1, $false, "26-14=12", 1) | Out-Null
$c = $t.Cell(6, 1)
$c.Range.Find.Execute("62-31=31", $true, $false, $false, $false, $false, $true, 1, $false, "31-23=8", 1) | Out-Null
$c = $t.Cell(6, 2)
$c.Range.Find.Execute("20-15=5", $true, $false, $false, $false, $false, $true, 1, $false, "42+43=85", 1) | Out-Null
$c = $t.Cell(6, 3)
$c.Range.Find.Execute("35-1=34", $true, $false, $false, $false, $false, $true, 1, $false, "76-39=37", 1) | Out-Null
$c = $t.Cell(6, 4)
$c.Range.Find.Execute("13-7=6", $true, $false, $false, $false, $false, $true, 1, $false, "25+64=89", 1) | Out-Null
$c = $t.Cell(6, 5)
$c.Range.Find.Execute("71-38=33", $true, $false, $false, $false, $false, $true, 1, $false, "69-64=5", 1) | Out-Null
$c = $t.Cell(7, 1)
$c.Range.Find.Execute("72-59=13", $true, $false, $false, $false, $false, $true, 1, $false, "61+20=81", 1) | Out-Null
$c = $t.Cell(7, 2)
$c.Range.Find.Execute("1+0=1", $true, $false, $false, $false, $false, $true, 1, $false, "42-34=8", 1) | Out-Null
$c = $t.Cell(7, 3)
$c.Range.Find.Execute("48-19=29", $true, $false, $false, $false, $false, $true, 1, $false, "14+44=58", 1) | Out-Null
$c = $t.Cell(7, 4)
$c.Range.Find.Execute("43+31=74", $true, $false, $false, $false, $false, $true, 1, $false, "65-26=39", 1) | Out-Null
$c = $t.Cell(7, 5)
$c.Range.Find.Execute("86-46=40", $true, $false, $false, $false, $false, $true, 1, $false, "60+9=69", 1) | Out-Null
$c = $t.Cell(8, 1)
$c.Range.Find.Execute("25+46=71", $true, $false, $false, $false, $false, $true, 1, $false, "15+46=61", 1) | Out-Null
$c = $t.Cell(8, 2)
$c.Range.Find.Execute("46-4=42", $true, $false, $false, $false, $false, $true, 1, $false, "4+65=69", 1) | Out-Null
$c = $t.Cell(8, 3)
$c.Range.Find.Execute("76-10=66", $true, $false, $false, $false, $false, $true, 1, $false, "59+7=66", 1) | Out-Null
$c = $t.Cell(8, 4)
$c.Range.Find.Execute("93-87=6", $true, $false, $false, $false, $false, $true, 1, $false, "85-80=5", 1) | Out-Null
$c = $t.Cell(8, 5)
$c.Range.Find.Execute("68-39=29", $true, $false, $false, $false, $false, $true, 1, $false, "9+82=91", 1) | Out-Null
$c = $t.Cell(9, 1)
$c.Range.Find.Execute("25+18=43", $true, $false, $false, $false, $false, $true, 1, $false, "54+26=80", 1) | Out-Null
$c = $t.Cell(9, 2)
$c.Range.Find.Execute("83-7=76", $true, $false, $false, $false, $false, $true, 1, $false, "22+10=32", 1) | Out-Null
$c = $t.Cell(9, 3)
$c.Range.Find.Execute("89-76=13", $true, $false, $false, $false, $false, $true, 1, $false, "77-54=23", 1) | Out-Null
$c = $t.Cell(9, 4)
$c.Range.Find.Execute("72-48=24", $true, $false, $false, $false, $false, $true, 1, $false, "40-21=19", 1) | Out-Null
$c = $t.Cell(9, 5)
$c.Range.Find.Execute("0+51=51", $true, $false, $false, $false, $false, $true, 1, $false, "1+90=91", 1) | Out-Null
$c = $t.Cell(10, 1)
$c.Range.Find.Execute("22+19=41", $true, $false, $false, $false, $false, $true, 1, $false, "58-45=13", 1) | Out-Null
$c = $t.Cell(10, 2)
$c.Range.Find.Execute("40-36=4", $true, $false, $false, $false, $false, $true, 1, $false, "77-25=52", 1) | Out-Null
$c = $t.Cell(10, 3)
$c.Range.Find.Execute("51-29=22", $true, $false, $false, $false, $false, $true, 1, $false, "82-81=1", 1) | Out-Null
$c = $t.Cell(10, 4)
$c.Range.Find.Execute("52+11=63", $true, $false, $false, $false, $false, $true, 1, $false, "73-48=25", 1) | Out-Null
$c = $t.Cell(10, 5)
$c.Range.Find.Execute("11+16=27", $true, $false, $false, $false, $false, $true, 1, $false, "85-62=23", 1) | Out-Null
$c = $t.Cell(11, 1)
$c.Range.Find.Execute("26+27=53", $true, $false, $false, $false, $false, $true, 1, $false, "41-26=15", 1) | Out-Null
$c = $t.Cell(11, 2)
$c.Range.Find.Execute("62+26=88", $true, $false, $false, $false, $false, $true, 1, $false, "1+53=54", 1) | Out-Null
$c = $t.Cell(11, 3)
$c.Range.Find.Execute("56+26=82", $true, $false, $false, $false, $false, $true, 1, $false, "6+84=90", 1) | Out-Null
$c = $t.Cell(11, 4)
$c.Range.Find.Execute("58-1=57", $true, $false, $false, $false, $false, $true, 1, $false, "47+38=85", 1) | Out-Null
$c = $t.Cell(11, 5)
$c.Range.Find.Execute("30+21=51", $true, $false, $false, $false, $false, $true, 1, $false, "20+29=49", 1) | Out-Null
$c = $t.Cell(12, 1)
$c.Range.Find.Execute("16+53=69", $true, $false, $false, $false, $false, $true, 1, $false, "43+34=77", 1) | Out-Null
$c = $t.Cell(12, 2)
$c.Range.Find.Execute("32+31=63", $true, $false, $false, $false, $false, $true, 1, $false, "96-79=17", 1) | Out-Null
$c = $t.Cell(12, 3)
$c.Range.Find.Execute("73-30=43", $true, $false, $false, $false, $false, $true, 1, $false, "97-50=47", 1) | Out-Null
$c = $t.Cell(12, 4)
$c.Range.Find.Execute("22+32=54", $true, $false, $false, $false, $false, $true, 1, $false, "14+27=41", 1) | Out-Null
$c = $t.Cell(12, 5)
$c.Range.Find.Execute("70+12=82", $true, $false, $false, $false, $false, $true, 1, $false, "18+71=89", 1) | Out-Null
$c = $t.Cell(13, 1)
$c.Range.Find.Execute("22+24=46", $true, $false, $false, $false, $false, $true, 1, $false, "28-3=25", 1) | Out-Null
$c = $t.Cell(13, 2)
$c.Range.Find.Execute("18+51=69", $true, $false, $false, $false, $false, $true, 1, $false, "86-71=15", 1) | Out-Null
$c = $t.Cell(13, 3)
$c.Range.Find.Execute("91-26=65", $true, $false, $false, $false, $false, $true, 1, $false, "9+20=29", 1) | Out-Null
$c = $t.Cell(13, 4)
$c.Range.Find.Execute("80+3=83", $true, $false, $false, $false, $false, $true, 1, $false, "33+55=88", 1) | Out-Null
$c = $t.Cell(13, 5)
$c.Range.Find.Execute("63+5=68", $true, $false, $false, $false, $false, $true, 1, $false, "31+14=45", 1) | Out-Null
$c = $t.Cell(14, 1)
$c.Range.Find.Execute("24+32=56", $true, $false, $false, $false, $false, $true, 1, $false, "16+56=72", 1) | Out-Null
$c = $t.Cell(14, 2)
$c.Range.Find.Execute("23+72=95", $true, $false, $false, $false, $false, $true, 1, $false, "28-11=17", 1) | Out-Null
$c = $t.Cell(14, 3)
$c.Range.Find.Execute("87-47=40", $true, $false, $false, $false, $false, $true, 1, $false, "10+82=92", 1) | Out-Null
$c = $t.Cell(14, 4)
$c.Range.Find.Execute("53+33=86", $true, $false, $false, $false, $false, $true, 1, $false, "89-65=24", 1) | Out-Null
$c = $t.Cell(14, 5)
$c.Range.Find.Execute("47+8=55", $true, $false, $false, $false, $false, $true, 1, $false, "12+83=95", 1) | Out-Null
$c = $t.Cell(15, 1)
$c.Range.Find.Execute("96-90=6", $true, $false, $false, $false, $false, $true, 1, $false, "26+46=72", 1) | Out-Null
$c = $t.Cell(15, 2)
$c.Range.Find.Execute("55-37=18", $true, $false, $false, $false, $false, $true, 1, $false, "3+48=51", 1) | Out-Null
$c = $t.Cell(15, 3)
$c.Range.Find.Execute("57-12=45", $true, $false, $false, $false, $false, $true, 1, $false, "9+89=98", 1) | Out-Null
$c = $t.Cell(15, 4)
$c.Range.Find.Execute("85-61=24", $true, $false, $false, $false, $false, $true, 1, $false, "77-3=74", 1) | Out-Null
$c = $t.Cell(15, 5)
$c.Range.Find.Execute("11+84=95", $true, $false, $false, $false, $false, $true, 1, $false, "89+10=99", 1) | Out-Null
$c = $t.Cell(16, 1)
$c.Range.Find.Execute("80-1=79", $true, $false, $false, $false, $false, $true, 1, $false, "81-50=31", 1) | Out-Null
$c = $t.Cell(16, 2)
$c.Range.Find.Execute("6+1=7", $true, $false, $false, $false, $false, $true, 1, $false, "65+10=75", 1) | Out-Null
$c = $t.Cell(16, 3)
$c.Range.Find.Execute("51-27=24", $true, $false, $false, $false, $false, $true, 1, $false, "59-12=47", 1) | Out-Null
$c = $t.Cell(16, 4)
$c.Range.Find.Execute("83-59=24", $true, $false, $false, $false, $false, $true, 1, $false, "20-14=6", 1) | Out-Null
$c = $t.Cell(16, 5)
$c.Range.Find.Execute("75-47=28", $true, $false, $false, $false, $false, $true, 1, $false, "53-5=48", 1) | Out-Null
$c = $t.Cell(17, 1)
$c.Range.Find.Execute("15+84=99", $true, $false, $false, $false, $false, $true, 1, $false, "69-0=69", 1) | Out-Null
$c = $t.Cell(17, 2)
$c.Range.Find.Execute("65-57=8", $true, $false, $false, $false, $false, $true, 1, $false, "15+9=24", 1) | Out-Null
$c = $t.Cell(17, 3)
$c.Range.Find.Execute("41+19=60", $true, $false, $false, $false, $false, $true, 1, $false, "52+40=92", 1) | Out-Null
$c = $t.Cell(17, 4)
$c.Range.Find.Execute("48-19=29", $true, $false, $false, $false, $false, $true, 1, $false, "52-45=7", 1) | Out-Null
$c = $t.Cell(17, 5)
$c.Range.Find.Execute("62-25=37", $true, $false, $false, $false, $false, $true, 1, $false, "89-57=32", 1) | Out-Null
$c = $t.Cell(18, 1)
$c.Range.Find.Execute("13+86=99", $true, $false, $false, $false, $false, $true, 1, $false, "9-5=4", 1) | Out-Null
$c = $t.Cell(18, 2)
$c.Range.Find.Execute("78-35=43", $true, $false, $false, $false, $false, $true, 1, $false, "37+20=57", 1) | Out-Null
$c = $t.Cell(18, 3)
$c.Range.Find.Execute("38-11=27", $true, $false, $false, $false, $false, $true, 1, $false, "8+57=65", 1) | Out-Null
$c = $t.Cell(18, 4)
$c.Range.Find.Execute("94-78=16", $true, $false, $false, $false, $false, $true, 1, $false, "52+22=74", 1) | Out-Null
$c = $t.Cell(18, 5)
$c.Range.Find.Execute("98-22=76", $true, $false, $false, $false, $false, $true, 1, $false, "80-45=35", 1) | Out-Null
$c = $t.Cell(19, 1)
$c.Range.Find.Execute("84-11=73", $true, $false, $false, $false, $false, $true, 1, $false, "53+18=71", 1) | Out-Null
$c = $t.Cell(19, 2)
$c.Range.Find.Execute("1+41=42", $true, $false, $false, $false, $false, $true, 1, $false, "92-22=70", 1) | Out-Null
$c = $t.Cell(19, 3)
$c.Range.Find.Execute("99-73=26", $true, $false, $false, $false, $false, $true, 1, $false, "93-37=56", 1) | Out-Null
$c = $t.Cell(19, 4)
$c.Range.Find.Execute("31-22=9", $true, $false, $false, $false, $false, $true, 1, $false, "63-32=31", 1) | Out-Null
$c = $t.Cell(19, 5)
$c.Range.Find.Execute("70-20=50", $true, $false, $false, $false, $false, $true, 1, $false, "15+60=75", 1) | Out-Null
$c = $t.Cell(20, 1)
$c.Range.Find.Execute("18+54=72", $true, $false, $false, $false, $false, $true, 1, $false, "57-2=55", 1) | Out-Null
$c = $t.Cell(20, 2)
$c.Range.Find.Execute("87+2=89", $true, $false, $false, $false, $false, $true, 1, $false, "39+56=95", 1) | Out-Null
$c = $t.Cell(20, 3)
$c.Range.Find.Execute("44-26=18", $true, $false, $false, $false, $false, $true, 1, $false, "35-16=19", 1) | Out-Null
$c = $t.Cell(20, 4)
$c.Range.Find.Execute("84-38=46", $true, $false, $false, $false, $false, $true, 1, $false, "56+16=72", 1) | Out-Null
$c = $t.Cell(20, 5)
$c.Range.Find.Execute("19+39=58", $true, $false, $false, $false, $false, $true, 1, $false, "60-41=19", 1) | Out-Null

Write-Host "done"
